$p = $ppt.ActivePresentation

$names = @("Parts", "OpenXmlPackage", "Package", "PackageParts", "XmlParts", "Document", "InnerXml", "OuterXml", "Xml")
foreach ($n in $names) {
    try {
        $v = $p.$n
        Write-Host "$n => $v"
    } catch {
        Write-Host "$n => ERROR $_"
    }
}
